# Apply the LinuxForHealth rebrand edit to the StructureDefinition workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-event-timestamp"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (Extension) - Constraint(s) text was removed (now duplicated only on row 4)
$elements.Range("AI2").Value = ""

# Row 5 (Extension.url) - Fixed Value URL updated to the new domain
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-event-timestamp"
